$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-observation row is inserted at row 44 ("Fruta / hortaliza, semanal"),
# pushing the former rows 44-57 down to become rows 45-58.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new observation.
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 44855
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 100112026
$ws.Range("G44").Value = "Haba"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 120
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 7500
$ws.Range("M44").Value = 7250
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Provincia de Diguillín"
$ws.Range("P44").Value = 290
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
